$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed "cryptos" price/volume snapshot (GitHub Actions bot run).
# Numeric-looking Price values are written with a leading "'" so Excel
# keeps them as text (matching the sheet's inlineStr/text-only columns)
# instead of auto-converting them to numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.945.11"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.501.34"
$ws.Range("E3").Value = "  -1.70%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'603.19"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'198.95"
$ws.Range("E6").Value = "  +6.83%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +1.23%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -3.10%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +1.42%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'54.38"
$ws.Range("E11").Value = "  +0.92%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "'0.0000303"
$ws.Range("E12").Value = "  -2.16%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'9.58"
$ws.Range("E13").Value = "  +0.47%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.059.62"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15 - BitcoinCash
$ws.Range("D15").Value = "'596.32"

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "69.955.92"
$ws.Range("E16").Value = "  -0.60%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "'19.01"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'12.65"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.479.25"
$ws.Range("E19").Value = "  -2.81%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "'17.98"
$ws.Range("E22").Value = "  +3.65%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'104.03"
$ws.Range("E23").Value = "  +10.35%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -2.02%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.53%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "'3.10"
$ws.Range("E26").Value = "  +5.35%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  +0.19%  "

# Row 28 - Filecoin
$ws.Range("D28").Value = "'9.84"
$ws.Range("E28").Value = "  +4.41%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +3.97%  "

# Row 30 - dogwifhat
$ws.Range("D30").Value = "'4.52"
$ws.Range("E30").Value = "  +23.06%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.30"
$ws.Range("E31").Value = "  +3.31%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "'12.76"
$ws.Range("E32").Value = "  +4.35%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.50%  "

# Row 34 - OKB
$ws.Range("D34").Value = "'63.61"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - Maker
$ws.Range("D35").Value = "3.693.32"
$ws.Range("E35").Value = "  +4.57%  "

# Row 36 - PEPE
$ws.Range("E36").Value = "  +2.66%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.14%  "

# Row 38 - Bittensor
$ws.Range("D38").Value = "'517.26"
$ws.Range("E38").Value = "  -2.17%  "

# Row 39 - TheGraph
$ws.Range("E39").Value = "  -2.91%  "

# Row 40/41 - InjectiveProtocol and Fetch.AI swapped places in the ranking
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  -6.47%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'36.74"
$ws.Range("E41").Value = "  -1.81%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -0.20%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  +0.79%  "

# Row 45 - ThetaToken
$ws.Range("E45").Value = "  -2.59%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  -0.70%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  -4.49%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "'8.79"
$ws.Range("E48").Value = "  -4.46%  "

# Row 49 - FirstDigitalUSD
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50 - Monero
$ws.Range("D50").Value = "'132.50"
$ws.Range("E50").Value = "  -3.10%  "

# Row 51 - FLOKI
$ws.Range("E51").Value = "  -1.62%  "
